$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original Text format so Excel does not
# auto-convert numeric-looking strings (e.g. "0.9934") into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 45 and 46 swap places: BabyDogeCoin moves up to rank 45, Aave moves down to rank 46.
# (rank/index column A stays the same, only Coin/Link/Price/Volume are updated)
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("E45").Value = "  +3.16%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "64.80"
$ws.Range("E46").Value = "  +0.12%  "

# Update Price (D) and Volume(1h) (E) columns for other changed rows
$ws.Range("D2").Value = "28.881.55"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.821.19"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").Value = "0.9934"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "242.82"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "0.6285"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.9945"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "0.07445"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "0.2928"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "0.07668"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.824.08"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "4.971"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "0.6651"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "82.75"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "0.000009664"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "6.006"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "28.916.54"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "224.88"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "0.9933"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "7.103"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").Value = "0.9947"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "159.84"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "0.1408"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").Value = "8.478"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "4.104"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "4.039"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "0.05444"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "1.849"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "0.7421"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Value = "2.603"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "1.236.03"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").Value = "2.735"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").Value = "0.01772"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "6.663"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").Value = "0.8970"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "0.9942"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "101.19"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "1.971.42"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D47").Value = "0.5059"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "0.4034"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").Value = "0.07406"
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("D50").Value = "8.917"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "1.655"
$ws.Range("E51").Value = "  +1.16%  "
